$d = $word.ActiveDocument

$old = "IC50=coeff1*ln<E(T)>"
$new = "<E(T)>=coeff1*ln(IC50)"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)

Write-Output "Found: $found"
